$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "28.175.80"
Set-TextValue $ws.Range("E2") "  +3.12%  "
Set-TextValue $ws.Range("D3") "1.777.22"
Set-TextValue $ws.Range("E3") "  -0.55%  "
Set-TextValue $ws.Range("D4") "1.003"
Set-TextValue $ws.Range("E4") "  -0.02%  "
Set-TextValue $ws.Range("D5") "339.08"
Set-TextValue $ws.Range("E5") "  -0.42%  "
Set-TextValue $ws.Range("D6") "1.001"
Set-TextValue $ws.Range("E6") "  +0.24%  "
Set-TextValue $ws.Range("D7") "0.3824"
Set-TextValue $ws.Range("E7") "  -3.21%  "
Set-TextValue $ws.Range("D8") "0.3424"
Set-TextValue $ws.Range("E8") "  -1.19%  "
Set-TextValue $ws.Range("D9") "47.00"
Set-TextValue $ws.Range("E9") "  -2.38%  "
Set-TextValue $ws.Range("D10") "1.146"
Set-TextValue $ws.Range("E10") "  -4.18%  "
Set-TextValue $ws.Range("D11") "0.07393"
Set-TextValue $ws.Range("E11") "  -1.34%  "
Set-TextValue $ws.Range("D12") "23.40"
Set-TextValue $ws.Range("E12") "  +7.30%  "
Set-TextValue $ws.Range("D13") "1.000"
Set-TextValue $ws.Range("E13") "  -0.05%  "
Set-TextValue $ws.Range("D14") "6.395"
Set-TextValue $ws.Range("E14") "  -1.35%  "
Set-TextValue $ws.Range("D15") "7.330"
Set-TextValue $ws.Range("E15") "  +3.02%  "
Set-TextValue $ws.Range("D16") "1.777.26"
Set-TextValue $ws.Range("E16") "  -0.37%  "
Set-TextValue $ws.Range("D17") "0.00001078"
Set-TextValue $ws.Range("E17") "  -1.62%  "
Set-TextValue $ws.Range("D18") "0.06658"
Set-TextValue $ws.Range("E18") "  -0.62%  "
Set-TextValue $ws.Range("D19") "82.66"
Set-TextValue $ws.Range("E19") "  -2.55%  "
Set-TextValue $ws.Range("D20") "1.002"
Set-TextValue $ws.Range("E20") "  +0.24%  "
Set-TextValue $ws.Range("D21") "17.40"
Set-TextValue $ws.Range("E21") "  -2.01%  "
Set-TextValue $ws.Range("D22") "6.407"
Set-TextValue $ws.Range("E22") "  -1.64%  "
Set-TextValue $ws.Range("D23") "28.162.94"
Set-TextValue $ws.Range("E23") "  +2.98%  "
Set-TextValue $ws.Range("D24") "12.11"
Set-TextValue $ws.Range("E24") "  -2.52%  "
Set-TextValue $ws.Range("D25") "2.378"
Set-TextValue $ws.Range("E25") "  -1.50%  "
Set-TextValue $ws.Range("D26") "20.73"
Set-TextValue $ws.Range("E26") "  -2.60%  "
Set-TextValue $ws.Range("D27") "1.430"
Set-TextValue $ws.Range("E27") "  -1.94%  "
Set-TextValue $ws.Range("D28") "2.408"
Set-TextValue $ws.Range("E28") "  -3.75%  "
Set-TextValue $ws.Range("D29") "154.06"
Set-TextValue $ws.Range("E29") "  -2.45%  "
Set-TextValue $ws.Range("D30") "1.978.29"
Set-TextValue $ws.Range("E30") "  -0.41%  "
Set-TextValue $ws.Range("D31") "134.62"
Set-TextValue $ws.Range("E31") "  -1.42%  "
Set-TextValue $ws.Range("D32") "4.016"
Set-TextValue $ws.Range("E32") "  -0.34%  "
Set-TextValue $ws.Range("D33") "6.070"
Set-TextValue $ws.Range("E33") "  +1.16%  "
Set-TextValue $ws.Range("D34") "0.08891"
Set-TextValue $ws.Range("E34") "  +0.67%  "
Set-TextValue $ws.Range("D35") "12.72"
Set-TextValue $ws.Range("E35") "  -2.42%  "
Set-TextValue $ws.Range("D36") "0.02414"
Set-TextValue $ws.Range("E36") "  -0.20%  "
Set-TextValue $ws.Range("D37") "0.6848"
Set-TextValue $ws.Range("E37") "  +0.24%  "
Set-TextValue $ws.Range("D38") "5.338"
Set-TextValue $ws.Range("E38") "  -1.41%  "
Set-TextValue $ws.Range("D39") "0.06363"
Set-TextValue $ws.Range("E39") "  -2.79%  "
Set-TextValue $ws.Range("D40") "0.2161"
Set-TextValue $ws.Range("E40") "  -2.42%  "
Set-TextValue $ws.Range("D41") "1.239"
Set-TextValue $ws.Range("E41") "  -1.18%  "
Set-TextValue $ws.Range("D42") "1.494"
Set-TextValue $ws.Range("E42") "  -7.63%  "
Set-TextValue $ws.Range("D43") "8.222"
Set-TextValue $ws.Range("E43") "  -1.90%  "
Set-TextValue $ws.Range("D44") "1.002"
Set-TextValue $ws.Range("E44") "  +0.28%  "
Set-TextValue $ws.Range("D45") "14.13"
Set-TextValue $ws.Range("E45") "  -2.57%  "
Set-TextValue $ws.Range("D46") "0.6270"
Set-TextValue $ws.Range("E46") "  -2.00%  "
Set-TextValue $ws.Range("D47") "3.862"
Set-TextValue $ws.Range("E47") "  -0.22%  "
Set-TextValue $ws.Range("D48") "132.50"
Set-TextValue $ws.Range("E48") "  -0.10%  "
Set-TextValue $ws.Range("D49") "2.068"
Set-TextValue $ws.Range("E49") "  -3.14%  "
Set-TextValue $ws.Range("D50") "0.07500"
Set-TextValue $ws.Range("E50") "  +4.70%  "
Set-TextValue $ws.Range("D51") "1.204"
Set-TextValue $ws.Range("E51") "  +3.74%  "
